$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh blank row above row 1. Every existing row (the old header
# row plus all 30 data rows) shifts down by one - row 1 -> row 2, ...,
# row 31 -> row 32 - carrying its original content/type untouched.
$ws.Rows.Item(1).Insert()

# The insert drags row 1's old bold/centered/thin-bordered header style down
# onto row 2 along with its content. Copy that exact formatting back up onto
# the new row 1 (format-only paste, so it reuses the existing style rather
# than inventing a near-duplicate one), then reset row 2 to the plain,
# unstyled look the data rows use.
$ws.Range("A2:J2").Copy()
$ws.Range("A1:J1").PasteSpecial(-4122)
$ws.Range("A2:J2").Style = "Normal"

# New row 1 becomes the numeric 0..9 header sequence.
for ($c = 1; $c -le 10; $c++) {
    $ws.Cells.Item(1, $c).Value = $c - 1
}

# Row 2 (old header row, now shifted down) keeps its text labels in columns
# A-F/H, but loses the values that used to live in G/I/J.
$ws.Cells.Item(2, 7).Value = $null
$ws.Cells.Item(2, 9).Value = $null
$ws.Cells.Item(2, 10).Value = $null
